$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the two rows of data ---
# Before: row1 = 1,2,3 (numbers)   row2 = one,two,three (strings)
# After:  row1 = one,two,three (strings)   row2 = 1,2,3 (numbers)
$ws.Range("A1").Value = "one"
$ws.Range("B1").Value = "two"
$ws.Range("C1").Value = "three"

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 3

# --- Touch the hyperlink styles (Hyperlink / Followed Hyperlink) get
# registered in the stylesheet, matching the look of a link having been
# inserted and then removed/un-styled again ---
$ws.Hyperlinks.Add($ws.Range("C1"), "http://example.com")
$ws.Hyperlinks(1).Delete()
$ws.Range("C1").ClearFormats()

# --- Selection moves to C1 ---
$ws.Range("C1").Select()

# --- Force a page setup block onto the sheet (paper size / orientation / dpi) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
